$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("gens")
$ws.Range("C10").Value = 100
$ws.Range("C24").Value = 0
$ws.Range("C33").Value = 0

$ws = $wb.Worksheets.Item("lines")
$ws.Range("C2").Value = 54.857333
$ws.Range("D2").Value = 0
$ws.Range("C3").Value = -153.85055
$ws.Range("D3").Value = 0
$ws.Range("C4").Value = -9.006787299999999
$ws.Range("D4").Value = 0
$ws.Range("C5").Value = -41.043154
$ws.Range("D5").Value = 0
$ws.Range("C6").Value = -1.0995133
$ws.Range("D6").Value = 0
$ws.Range("C7").Value = 134.90336
$ws.Range("D7").Value = 0
$ws.Range("C8").Value = -24.753901
$ws.Range("D8").Value = 0
$ws.Range("C9").Value = -115.04315
$ws.Range("D9").Value = 0
$ws.Range("C10").Value = -80.006787
$ws.Range("D10").Value = 0
$ws.Range("C11").Value = -137.09951
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("C12").Value = 175
$ws.Range("F12").Value = 0
$ws.Range("C13").Value = -24.069065
$ws.Range("D13").Value = 0
$ws.Range("C14").Value = 28.069065
$ws.Range("D14").Value = 0
$ws.Range("C15").Value = -70.809938
$ws.Range("D15").Value = 0
$ws.Range("C16").Value = -108.39893
$ws.Range("D16").Value = 0
$ws.Range("C17").Value = -173.22412
$ws.Range("D17").Value = 0.00000000000010450204
$ws.Range("C18").Value = -210.81311
$ws.Range("D18").Value = -0.00000000000010450204
$ws.Range("C19").Value = -233.81409
$ws.Range("D19").Value = 0.000000000000059715453
$ws.Range("C20").Value = -10.219968
$ws.Range("D20").Value = 0
$ws.Range("C21").Value = -168.03337
$ws.Range("D21").Value = -0.000000000000059715453
$ws.Range("C22").Value = -151.17867
$ws.Range("D22").Value = 0
$ws.Range("C23").Value = -75.84746199999999
$ws.Range("D23").Value = 0
$ws.Range("C24").Value = -204.21997
$ws.Range("D24").Value = 0
$ws.Range("C25").Value = 58.457549
$ws.Range("D25").Value = 0
$ws.Range("C26").Value = -122.60573
$ws.Range("D26").Value = 0.000000000000014210855
$ws.Range("C27").Value = -122.60573
$ws.Range("D27").Value = -0.000000000000014210855
$ws.Range("C28").Value = 24.753901
$ws.Range("D28").Value = 0
$ws.Range("C29").Value = -121.78855
$ws.Range("D29").Value = 0
$ws.Range("C30").Value = 31.02613
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("C31").Value = 18.297905
$ws.Range("D31").Value = 0.0000000000000011500114
$ws.Range("E31").Value = 0
$ws.Range("C32").Value = -140.08645
$ws.Range("D32").Value = -0.0000000000000086250852
$ws.Range("C33").Value = -157.35105
$ws.Range("D33").Value = 0.0000000000000021357354
$ws.Range("C34").Value = -157.35105
$ws.Range("D34").Value = 0
$ws.Range("C35").Value = -74.986935
$ws.Range("D35").Value = 0
$ws.Range("C36").Value = -74.986935
$ws.Range("D36").Value = 0
$ws.Range("C37").Value = -138.98693
$ws.Range("D37").Value = 0
$ws.Range("C38").Value = -138.98693
$ws.Range("D38").Value = 0
$ws.Range("C39").Value = -159.91355
$ws.Range("D39").Value = 0.0000000000000055857695

$ws = $wb.Worksheets.Item("bus")
$ws.Range("B2").Value = 115.09
$ws.Range("C2").Value = -124.83327
$ws.Range("B3").Value = 115.09
$ws.Range("C3").Value = -125.60127
$ws.Range("B4").Value = 115.09
$ws.Range("C4").Value = -92.370802
$ws.Range("B5").Value = 115.09
$ws.Range("C5").Value = -120.38879
$ws.Range("B6").Value = 115.09
$ws.Range("C6").Value = -124.06769
$ws.Range("B7").Value = 115.09
$ws.Range("C7").Value = -125.39016
$ws.Range("B8").Value = 115.09
$ws.Range("C8").Value = -101.7207
$ws.Range("B9").Value = 115.09
$ws.Range("C9").Value = -112.3957
$ws.Range("B10").Value = 115.09
$ws.Range("C10").Value = -108.4243
$ws.Range("B11").Value = 115.09
$ws.Range("C11").Value = -117.02709
$ws.Range("B12").Value = 115.09
$ws.Range("C12").Value = -102.47627
$ws.Range("B13").Value = 115.09
$ws.Range("C13").Value = -99.318792
$ws.Range("B14").Value = 115.09
$ws.Range("C14").Value = -91.25319
$ws.Range("B15").Value = 115.09
$ws.Range("C15").Value = -102.04703
$ws.Range("B16").Value = 115.09
$ws.Range("C16").Value = -89.004271
$ws.Range("B17").Value = 115.09
$ws.Range("C17").Value = -89.99805000000001
$ws.Range("B18").Value = 115.09
$ws.Range("C18").Value = -86.831547
$ws.Range("B19").Value = 115.09
$ws.Range("C19").Value = -87.087718
$ws.Range("B20").Value = 115.09
$ws.Range("C20").Value = -90.711651
$ws.Range("B21").Value = 115.09
$ws.Range("C21").Value = -87.71217300000001
$ws.Range("B22").Value = 115.09
$ws.Range("C22").Value = -82.996591
$ws.Range("B23").Value = 115.09
$ws.Range("C23").Value = -72.12247000000001
$ws.Range("B24").Value = 115.09
$ws.Range("C24").Value = -84.654461
$ws.Range("B25").Value = 115.09
$ws.Range("C25").Value = -90.29147399999999

